$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Neg_Change")
$data1 = @(
    ,@("HEROMOTOCO", 5099.5, 5148.3, 5061.4, 5109.5, 783364, 1584393, -0.505574690117919, "HEROMOTOCO")
    ,@("BPCL", 310.7, 312.3, 307, 308.55, 4334416, 8888325, -0.5123472645295936, "BPCL")
    ,@("ICICIGI", 1852, 1868.3, 1832.1, 1843, 441249, 870926, -0.4933564964187543, "ICICIGI")
    ,@("PATANJALI", 1780, 1796.4, 1771.3, 1777.6, 230722, 571455, -0.5962551731982396, "PATANJALI")
    ,@("FORTIS", 911.4, 920.85, 904.4, 915, 963721, 1982246, -0.5138237131011993, "FORTIS")
    ,@("PAGEIND", 44275, 44575, 43650, 44120, 21171, 47648, -0.5556791470785762, "PAGEIND")
    ,@("CONCOR", 531.25, 533.9, 526.25, 526.9, 826127, 2055711, -0.5981307683813533, "CONCOR")
    ,@("MAXHEALTH", 1151, 1173.9, 1147.6, 1150, 1977662, 4412704, -0.551825366034069, "MAXHEALTH")
    ,@("GLENMARK", 1936, 1947.8, 1910.7, 1919.7, 239914, 481992, -0.5022448505369383, "GLENMARK")
    ,@("KALYANKJIL", 508, 510.5, 501.5, 503, 3445371, 7354852, -0.5315512807055804, "KALYANKJIL")
    ,@("LAURUSLABS", 864, 872.6, 854.75, 855, 646790, 1373559, -0.5291137839728763, "LAURUSLABS")
    ,@("PNBHOUSING", 770.1, 774.5, 750.3, 753.9, 1227892, 2816692, -0.564065932661434, "PNBHOUSING")
)
for ($i = 0; $i -lt $data1.Count; $i++) {
    $r = $i + 2
    $row = $data1[$i]
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $ws1.Cells.Item($r, 9).Value = $row[8]
}

$ws2 = $wb.Worksheets.Item("Pos_Change")
$data2 = @(
    ,@("ETERNAL", 315.55, 316.75, 307.75, 315.6, 62441100, 44020764, 0.4184465312778306, "ETERNAL")
    ,@("NESTLEIND", 1160.3, 1184.4, 1150, 1154.2, 1593255, 1065888, 0.4947677429523552, "NESTLEIND")
    ,@("INFY", 1501.9, 1501.9, 1467.7, 1469, 8169075, 5443364, 0.5007401672936074, "INFY")
    ,@("RELIANCE", 1381.1, 1403.5, 1350, 1356, 18758842, 13161813, 0.4252475703765127, "RELIANCE")
    ,@("ADANIENSOL", 775.4, 776.4, 759.25, 764, 1341814, 891046, 0.5058863403236197, "ADANIENSOL")
    ,@("JINDALSTEL", 970.4, 974.4, 943.55, 949.1, 968556, 673705, 0.4376559473360002, "JINDALSTEL")
    ,@("NAUKRI", 1390, 1395, 1354.5, 1356.5, 1215734, 830602, 0.4636781515093872, "NAUKRI")
    ,@("PERSISTENT", 5308.5, 5329.5, 5228, 5320, 487634, 316243, 0.5419598220355866, "PERSISTENT")
    ,@("SRF", 2832, 2869.9, 2828.6, 2852, 537875, 363957, 0.4778531529823578, "SRF")
    ,@("MANKIND", 2488, 2511.3, 2455.1, 2474.3, 330123, 233219, 0.415506455305957, "MANKIND")
    ,@("APLAPOLLO", 1599.1, 1633.6, 1578, 1600, 720090, 509123, 0.4143733439659965, "APLAPOLLO")
    ,@("MCX", 7387, 7463, 7304, 7380, 644963, 455408, 0.416231159751256, "MCX")
    ,@("IIFL", 436, 440.45, 425.9, 427.1, 1027543, 650392, 0.5798825938818436, "IIFL")
)
for ($i = 0; $i -lt $data2.Count; $i++) {
    $r = $i + 2
    $row = $data2[$i]
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $ws2.Cells.Item($r, 6).Value = $row[5]
    $ws2.Cells.Item($r, 7).Value = $row[6]
    $ws2.Cells.Item($r, 8).Value = $row[7]
    $ws2.Cells.Item($r, 9).Value = $row[8]
}
